$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new row (each / helps iterating through an array) ---------------
$ws.Range("C8").Value = 'each'
$ws.Range("B5").Value = 'Search from list of similar locators'
$ws.Range("B6").Value = 'Assertion'
$ws.Range("B7").Value = 'Helps finding some specific locator'
$ws.Range("B8").Value = 'Helps iterating through an array'

# --- E14: re-quote ---------------------------------------------------------
$ws.Range("E14").Value = '"should" is the assertion type of "Chai"'

# --- Row 9: new row (wrap / cy.warp) ---------------------------------------
$ws.Range("B9").Value = 'To resolve the promise, thus refraining from getting ''click()'' method deprecated'
$ws.Range("D9").Value = 'cy.warp(    )'
$ws.Range("C9").Value = 'wrap'

# --- E16 / E18 notes ---------------------------------------------------------
$ws.Range("E16").Value = 'Cypress is asynchronous in nature and there is no guarantee in sequence of execution , but Cypress takes care of it.'
$ws.Range("E18").Value = 'Promise comes with ''resolved'', ''rejected'' and ''pending'''

# --- Row 10: new row (then / .then()) ---------------------------------------
$ws.Range("C10").Value = 'then'
$ws.Range("D10").Value = '.then()'
$ws.Range("B10").Value = 'Wait until promise is resolved(don''t rush)'

# --- Row 11: new row (text / .text()) ----------------------------------------
$ws.Range("C11").Value = 'text'
$ws.Range("D11").Value = '.text()'
$ws.Range("B11").Value = 'jQuery command return text content of the selected elements(supported by Cypress after manually resolving promise)'

# --- E20 note -----------------------------------------------------------------
$ws.Range("E20").Value = 'Non cypress commands can not resolve promise by themselves, we need to manually resolve it by using then()'

# --- Row heights for the wrapped, multi-line rows ------------------------------
$ws.Rows.Item(9).RowHeight = 28.8
$ws.Rows.Item(11).RowHeight = 43.2

# --- Final selection matches the authored state --------------------------------
$ws.Range("E21").Select()
